$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.397.80'
$ws.Range("E2").Value = '  -3.03%  '
$ws.Range("D3").Value = '1.747.46'
$ws.Range("E3").Value = '  -3.44%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'322.27"
$ws.Range("E5").Value = '  -2.03%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = "'0.4238"
$ws.Range("E7").Value = '  -4.98%  '
$ws.Range("D8").Value = "'0.3597"
$ws.Range("E8").Value = '  -2.98%  '
$ws.Range("D9").Value = "'0.07492"
$ws.Range("E9").Value = '  -2.79%  '
$ws.Range("D10").Value = "'42.03"
$ws.Range("E10").Value = '  -6.44%  '
$ws.Range("D11").Value = "'1.098"
$ws.Range("E11").Value = '  -2.55%  '
$ws.Range("D13").Value = "'20.64"
$ws.Range("E13").Value = '  -6.42%  '
$ws.Range("D14").Value = "'6.024"
$ws.Range("E14").Value = '  -4.35%  '
$ws.Range("D15").Value = "'7.214"
$ws.Range("E15").Value = '  -4.71%  '
$ws.Range("D16").Value = '1.748.12'
$ws.Range("E16").Value = '  -5.48%  '
$ws.Range("D17").Value = "'92.82"
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("D19").Value = "'0.06365"
$ws.Range("E19").Value = '  -2.83%  '
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = "'17.03"
$ws.Range("E21").Value = '  -2.72%  '
$ws.Range("D22").Value = "'5.880"
$ws.Range("E22").Value = '  -5.57%  '
$ws.Range("D23").Value = '27.466.81'
$ws.Range("E23").Value = '  -2.98%  '
$ws.Range("E24").Value = '  -4.20%  '
$ws.Range("D25").Value = "'2.090"
$ws.Range("E25").Value = '  -3.46%  '
$ws.Range("D26").Value = "'161.97"
$ws.Range("E26").Value = '  +3.85%  '
$ws.Range("D27").Value = "'20.28"
$ws.Range("E27").Value = '  -2.46%  '
$ws.Range("D28").Value = '1.948.80'
$ws.Range("E28").Value = '  -4.58%  '
$ws.Range("D29").Value = "'2.126"
$ws.Range("E29").Value = '  -8.02%  '
$ws.Range("D30").Value = "'123.83"
$ws.Range("E30").Value = '  -3.53%  '
$ws.Range("D31").Value = "'1.102"
$ws.Range("E31").Value = '  -8.05%  '
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("D33").Value = "'5.528"
$ws.Range("E33").Value = '  -6.41%  '
$ws.Range("D34").Value = "'0.08896"
$ws.Range("E34").Value = '  -3.87%  '
$ws.Range("D35").Value = "'12.19"
$ws.Range("E35").Value = '  -6.49%  '
$ws.Range("D36").Value = "'0.02281"
$ws.Range("E36").Value = '  -3.04%  '
$ws.Range("E37").Value = '  -3.99%  '
$ws.Range("D38").Value = "'0.05994"
$ws.Range("E38").Value = '  -3.75%  '
$ws.Range("D39").Value = "'0.6329"
$ws.Range("E39").Value = '  -3.73%  '
$ws.Range("D40").Value = "'4.944"
$ws.Range("E40").Value = '  -4.57%  '
$ws.Range("D41").Value = "'1.182"
$ws.Range("E41").Value = '  -1.43%  '
$ws.Range("D42").Value = "'0.9999"
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = "'7.866"
$ws.Range("E43").Value = '  -3.35%  '
$ws.Range("D44").Value = "'1.384"
$ws.Range("E44").Value = '  -1.55%  '
$ws.Range("D45").Value = "'13.38"
$ws.Range("E45").Value = '  -4.08%  '
$ws.Range("D46").Value = "'0.5869"
$ws.Range("E46").Value = '  -3.57%  '
$ws.Range("E47").Value = '  -2.00%  '
$ws.Range("D48").Value = "'1.968"
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").Value = "'122.28"
$ws.Range("E49").Value = '  -3.76%  '
$ws.Range("D50").Value = "'1.167"
$ws.Range("E50").Value = '  +0.85%  '
$ws.Range("D51").Value = "'0.06803"
$ws.Range("E51").Value = '  -2.67%  '
